$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hospitalnotification")

# Update the two date values in column C (serial date 44223 -> 44224)
$ws.Range("C2").Value = 44224
$ws.Range("C3").Value = 44224

# Update the active selection cell on the sheet
$ws.Range("D14").Select()

# Update the workbook window position (yWindow 4160 -> 11420)
$excel.Windows.Item(1).Top = 11420
